# Updated cryptos list (Price / Volume(1h) columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values look like plain numbers (e.g. "7.20", "0.999") and
# Excel would silently coerce them to numeric cells, dropping the exact
# text (trailing zeros etc.). Mark those target cells as Text first so the
# literal string is preserved, matching the original inline-string cells.
$textCells = @(
    "D5", "D6", "D9", "D10", "D14", "D15", "D19", "D21", "D22", "D24", "D28", "D30", "D31", "D34", "D35", "D36", "D37", "D38", "D40", "D42", "D43", "D47", "D49", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values row by row.

$ws.Range('D2').Value = '65.657.91'
$ws.Range('E2').Value = '  -1.16%  '

$ws.Range('D3').Value = '3.451.93'
$ws.Range('E3').Value = '  -3.75%  '
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '595.92'
$ws.Range('E5').Value = '  -1.57%  '

$ws.Range('D6').Value = '137.06'
$ws.Range('E6').Value = '  -7.29%  '

$ws.Range('D7').Value = '3.454.17'
$ws.Range('E7').Value = '  -3.66%  '
$ws.Range('E8').Value = '  -0.10%  '

$ws.Range('D9').Value = '0.493'
$ws.Range('E9').Value = '  +0.23%  '

$ws.Range('D10').Value = '7.49'
$ws.Range('E10').Value = '  -5.62%  '
$ws.Range('E11').Value = '  -9.64%  '
$ws.Range('E12').Value = '  -7.66%  '

$ws.Range('D13').Value = '4.041.04'
$ws.Range('E13').Value = '  -3.70%  '

$ws.Range('D14').Value = '0.0000183'
$ws.Range('E14').Value = '  -10.35%  '

$ws.Range('D15').Value = '26.71'
$ws.Range('E15').Value = '  -9.53%  '

$ws.Range('D16').Value = '3.440.50'
$ws.Range('E16').Value = '  -4.43%  '

$ws.Range('D17').Value = '65.655.54'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('E18').Value = '  -2.18%  '

$ws.Range('D19').Value = '9.92'
$ws.Range('E19').Value = '  -10.03%  '
$ws.Range('E20').Value = '  -8.37%  '

$ws.Range('D21').Value = '13.77'
$ws.Range('E21').Value = '  -7.12%  '

$ws.Range('D22').Value = '396.26'
$ws.Range('E22').Value = '  -6.44%  '
$ws.Range('E23').Value = '  -10.01%  '

$ws.Range('D24').Value = '73.54'
$ws.Range('E24').Value = '  -5.87%  '
$ws.Range('E25').Value = '  -0.02%  '

$ws.Range('D26').Value = '3.602.25'
$ws.Range('E27').Value = '  -10.20%  '

$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('E29').Value = '  -8.70%  '

$ws.Range('D30').Value = '7.20'
$ws.Range('E30').Value = '  -10.69%  '

$ws.Range('D31').Value = '8.26'
$ws.Range('E31').Value = '  -11.62%  '

$ws.Range('D32').Value = '3.460.68'
$ws.Range('E32').Value = '  -3.50%  '

$ws.Range('D34').Value = '0.146'
$ws.Range('E34').Value = '  -6.44%  '

$ws.Range('D35').Value = '22.96'
$ws.Range('E35').Value = '  -7.89%  '

$ws.Range('D36').Value = '173.14'
$ws.Range('E36').Value = '  -0.93%  '

$ws.Range('D37').Value = '1.22'
$ws.Range('E37').Value = '  -13.41%  '

$ws.Range('D38').Value = '6.93'
$ws.Range('E38').Value = '  -10.40%  '
$ws.Range('E39').Value = '  -7.14%  '

$ws.Range('D40').Value = '4.83'
$ws.Range('E40').Value = '  -12.42%  '
$ws.Range('E41').Value = '  -8.13%  '

$ws.Range('D42').Value = '0.822'
$ws.Range('E42').Value = '  -6.45%  '

$ws.Range('D43').Value = '43.54'
$ws.Range('E43').Value = '  -5.50%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  -13.94%  '
$ws.Range('E46').Value = '  -11.06%  '

$ws.Range('D47').Value = '23.31'
$ws.Range('E47').Value = '  -2.56%  '
$ws.Range('E48').Value = '  -2.96%  '

$ws.Range('D49').Value = '6.58'
$ws.Range('E49').Value = '  -7.54%  '

$ws.Range('D50').Value = '2.12'
$ws.Range('E50').Value = '  -15.54%  '

$ws.Range('D51').Value = '2.222.57'
$ws.Range('E51').Value = '  -7.71%  '

# Drop back to the default cell style so the forced text format above
# does not leave a visible style change on these cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
